# Insert a new data row into the "Mango" price sheet.
# Before the edit, row 366 held a record that (together with every row
# below it up to 447) needs to shift down by one position so that a brand
# new record can be placed at row 366.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 366 (and everything after it) down by one row.
$ws.Rows("366").Insert()

# The record that used to be on row 366 is now on row 367; re-use its
# unchanged fields (A,B,C,E,F,G,H,I,J,K,L,Q,T) as the template for the new
# row 366, then overwrite the fields that actually differ for the new
# record.
$ws.Range("A366:T366").Value2 = $ws.Range("A367:T367").Value2
$ws.Range("D366").NumberFormat = $ws.Range("D367").NumberFormat

$ws.Range("D366").Value2 = 45204
$ws.Range("M366").Value2 = 100
$ws.Range("N366").Value2 = 12000
$ws.Range("O366").Value2 = 12000
$ws.Range("P366").Value2 = 12000
$ws.Range("R366").Value2 = "Brasil"
$ws.Range("S366").Value2 = 3000
